# Petty cash book 2021 - 16-Apr-2021 midday update.
# Adds new petty-cash entries for 44300 (16-Apr), 44301 (17-Apr) and
# 44302 (18-Apr) on Sheet1, and moves the frozen-pane / selection to the
# newly active working area.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---- 44300 (row 18) continuation -------------------------------------
$ws.Range("B18").Value = "Wages Expense"
$ws.Range("D18").Formula = "=60000+260000"

$ws.Range("B19").Value = "A/R"
$ws.Range("C19").Formula = "=1170000+3615000+6235000+6777500"

$ws.Range("B20").Value = "TRANSFER BCA"
$ws.Range("D20").Formula = "=1170000+3615000+140000+480500"

$ws.Range("B21").Value = "BELI lampu"
$ws.Range("D21").Value = 105000

$ws.Range("B22").Value = "SALES - cash/retail"
$ws.Range("C22").Formula = "=21068025-6446525-6777500"

$ws.Range("B23").Value = "FREIGHT OUT"
$ws.Range("D23").Formula = "=150000"

$ws.Range("B24").Value = "SELISIH - lebih"
$ws.Range("C24").Value = 10000

$ws.Range("B25").Value = "SETOR KE BANK"
$ws.Range("D25").Value = 20000000

# ---- 44301 (row 26) ----------------------------------------------------
$ws.Range("A26").Value = 44301
$ws.Range("B26").Value = "Wages Expense"
$ws.Range("D26").Formula = "=60000+280000"

$ws.Range("B27").Value = "TRANSFER BCA"
$ws.Range("D27").Formula = "=13015000+448000+300000000+4800000+17000000+31001000+5400000"

$ws.Range("B28").Value = "A/R"
$ws.Range("C28").Formula = "=300000000+31001000+42923000"

$ws.Range("B29").Value = "TELPON - 5224823"
$ws.Range("D29").Value = 330000

$ws.Range("B30").Value = "PLN - Astar 165"
$ws.Range("D30").Value = 817500

$ws.Range("B31").Value = "PLN - Astar 214"
$ws.Range("D31").Value = 103000

$ws.Range("B32").Value = "SALES - cash/retail"
$ws.Range("C32").Formula = "=41305475+6108025-42923000"

$ws.Range("B33").Value = "SELISIH - lebih"
$ws.Range("C33").Value = 45000

$ws.Range("B34").Value = "SETOR KE BANK"
$ws.Range("D34").Value = 5000000

# ---- 44302 (row 35) ----------------------------------------------------
$ws.Range("A35").Value = 44302
$ws.Range("B35").Value = "Wages Expense"
$ws.Range("D35").Formula = "=60000"

$ws.Range("B36").Value = "TRANSFER BCA"
$ws.Range("D36").Formula = "=5000000+9750000"

$ws.Range("B37").Value = "A/R"
$ws.Range("C37").Formula = "=9750000"

$ws.Range("B38").Value = "FREIGHT OUT"
$ws.Range("D38").Formula = "=63000"

# ---- view state: scroll working area into the frozen pane -------------
$ws.Activate()
$ws.Range("B39").Select()
$excel.ActiveWindow.ScrollRow = 34
